# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by copying the "2021-Q4" sheet
#    (same column layout/styles) and placing it right before "总计".
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$totalIndex = $total.Index
$template.Copy($total)

# NOTE: adding/copying a sheet shifts worksheet indices, so any handle
# obtained before the copy (like $total) may now resolve to the wrong
# sheet. The copy lands in the slot "总计" used to occupy (immediately
# before it), so grab it positionally and re-fetch everything else by
# name afterwards.
$newSheet = $wb.Worksheets.Item($totalIndex)
$newSheet.Name = "2022-Q1"

# Update the fund-holding data row on the new sheet. The numeric-looking
# fields (fund scale / position / market value) are stored as TEXT in
# this workbook's convention, so force text entry with a leading
# apostrophe and then strip the auto-applied "quote prefix" style back
# off so the cell matches the plain (unstyled) neighbouring cells.
$newSheet.Range("D2").Value = "'30.93"
$newSheet.Range("E2").Value = "'89.81"
$newSheet.Range("F2").Value = "'4.46"
$newSheet.Range("G2").Value = "'1.3795"
$newSheet.Range("D2:G2").ClearFormats()
$newSheet.Range("H2").Value = 8

# ------------------------------------------------------------------
# 2. Update the "总计" (total) sheet: insert a new row for 2022-Q1 at
#    the top of the data and keep the older quarters below it.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()

# Fix up the running index in column A (0,1,2,3,4,5).
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# Re-apply the header-like formatting to the new A2 cell.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# New row data, then clear the borrowed formatting off the new row so
# it matches the plain (unstyled) data rows below it.
$total.Range("B2:D2").ClearFormats()
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 1.38
